$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right before the current row 852 (shifts the existing
# rows 852:894 down to 854:896, matching the target dimension A1:R896).
$ws.Rows.Item(852).Insert()
$ws.Rows.Item(852).Insert()

# New weekly price entries (Escarola, 2022-07-04 / serial 44746).
$newRow1 = @(1, "Agrícola del Norte S.A. de Arica", "Arica y Parinacota", 44746, 15, 100112033, "Lechuga", "Escarola", "Primera", 120, 8000, 9000, 8500, "`$/caja 12 unidades", "Región de Arica y Parinacota", 708, 12, "Hortaliza")
$newRow2 = @(1, "Agrícola del Norte S.A. de Arica", "Arica y Parinacota", 44746, 15, 100112033, "Lechuga", "Escarola", "Segunda", 120, 8000, 9000, 8500, "`$/caja 18 unidades", "Región de Arica y Parinacota", 472, 18, "Hortaliza")

for ($i = 0; $i -lt $newRow1.Length; $i++) {
    $ws.Cells.Item(852, $i + 1).Value = $newRow1[$i]
}
for ($i = 0; $i -lt $newRow2.Length; $i++) {
    $ws.Cells.Item(853, $i + 1).Value = $newRow2[$i]
}
